# Auto-generated Excel COM-interop edit script
# Implements: actualizacion automatica del tracker (tracker_resultados)
#  1) Fill G120/H120 and G121/H121 results (Fallo / -1)
#  2) Append new match rows 129-150, extending the sheet dimension to A1:H150

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$text)
    # Force literal text even for values that look like dates ("2025-09-05", etc.)
    # so Excel does not silently reinterpret them as date serials, then strip the
    # helper number format back off so the cell is left with no special style.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# --- 1) Backfill results for two previously-pending rows ---
Set-TextCell $ws.Cells.Item(120, 7) "Fallo"
$ws.Cells.Item(120, 8).Value = -1

Set-TextCell $ws.Cells.Item(121, 7) "Fallo"
$ws.Cells.Item(121, 8).Value = -1

# --- 2) Append new rows 129-150 ---
# row 129
$ws.Cells.Item(129, 1).Value = 14494996
Set-TextCell $ws.Cells.Item(129, 2) "2025-09-05"
Set-TextCell $ws.Cells.Item(129, 3) "Naomi Osaka"
Set-TextCell $ws.Cells.Item(129, 4) "Amanda Anisimova"
Set-TextCell $ws.Cells.Item(129, 5) "Gana Amanda Anisimova"
$ws.Cells.Item(129, 6).Value = 2

# row 130
$ws.Cells.Item(130, 1).Value = 14552520
Set-TextCell $ws.Cells.Item(130, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(130, 3) "Oliver Crawford"
Set-TextCell $ws.Cells.Item(130, 4) "Rio Noguchi"
Set-TextCell $ws.Cells.Item(130, 5) "Gana Rio Noguchi"
$ws.Cells.Item(130, 6).Value = 2.5

# row 131
$ws.Cells.Item(131, 1).Value = 14552624
Set-TextCell $ws.Cells.Item(131, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(131, 3) "Kai Wehnelt"
Set-TextCell $ws.Cells.Item(131, 4) "Francesco Maestrelli"
Set-TextCell $ws.Cells.Item(131, 5) "Gana Kai Wehnelt"
$ws.Cells.Item(131, 6).Value = 5

# row 132
$ws.Cells.Item(132, 1).Value = 14609410
Set-TextCell $ws.Cells.Item(132, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(132, 3) "Pietro Romeo Scomparin"
Set-TextCell $ws.Cells.Item(132, 4) "Alessandro Coccioli"
Set-TextCell $ws.Cells.Item(132, 5) "Gana Alessandro Coccioli"
$ws.Cells.Item(132, 6).Value = 3.75

# row 133
$ws.Cells.Item(133, 1).Value = 14609413
Set-TextCell $ws.Cells.Item(133, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(133, 3) "Daniele Rapagnetta"
Set-TextCell $ws.Cells.Item(133, 4) "Michele Ribecai"
Set-TextCell $ws.Cells.Item(133, 5) "Gana Daniele Rapagnetta"
$ws.Cells.Item(133, 6).Value = 5

# row 134
$ws.Cells.Item(134, 1).Value = 14609412
Set-TextCell $ws.Cells.Item(134, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(134, 3) "Lorenzo Comino"
Set-TextCell $ws.Cells.Item(134, 4) "Juan Cruz Martin Manzano"
Set-TextCell $ws.Cells.Item(134, 5) "Gana Lorenzo Comino"
$ws.Cells.Item(134, 6).Value = 4.33

# row 135
$ws.Cells.Item(135, 1).Value = 14609424
Set-TextCell $ws.Cells.Item(135, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(135, 3) "Stefano Napolitano"
Set-TextCell $ws.Cells.Item(135, 4) "Ainius Sabaliauskas"
Set-TextCell $ws.Cells.Item(135, 5) "Gana Ainius Sabaliauskas"
$ws.Cells.Item(135, 6).Value = 4.33

# row 136
$ws.Cells.Item(136, 1).Value = 14609554
Set-TextCell $ws.Cells.Item(136, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(136, 3) "Stefan Palosi"
Set-TextCell $ws.Cells.Item(136, 4) "Jeremy Gschwendtner"
Set-TextCell $ws.Cells.Item(136, 5) "Gana Jeremy Gschwendtner"
$ws.Cells.Item(136, 6).Value = 3.75

# row 137
$ws.Cells.Item(137, 1).Value = 14609559
Set-TextCell $ws.Cells.Item(137, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(137, 3) "Edouard Villoslada"
Set-TextCell $ws.Cells.Item(137, 4) "Mihai Alexandru Coman"
Set-TextCell $ws.Cells.Item(137, 5) "Gana Mihai Alexandru Coman"
$ws.Cells.Item(137, 6).Value = 2.25

# row 138
$ws.Cells.Item(138, 1).Value = 14609558
Set-TextCell $ws.Cells.Item(138, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(138, 3) "Matei Florin Breazu"
Set-TextCell $ws.Cells.Item(138, 4) "Stefan Adrian Andreescu"
Set-TextCell $ws.Cells.Item(138, 5) "Gana Matei Florin Breazu"
$ws.Cells.Item(138, 6).Value = 6.5

# row 139
$ws.Cells.Item(139, 1).Value = 14610355
Set-TextCell $ws.Cells.Item(139, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(139, 3) "Adhithya Ganesan"
Set-TextCell $ws.Cells.Item(139, 4) "Leyton Rivera"
Set-TextCell $ws.Cells.Item(139, 5) "Gana Leyton Rivera"
$ws.Cells.Item(139, 6).Value = 3.75

# row 140
$ws.Cells.Item(140, 1).Value = 14610354
Set-TextCell $ws.Cells.Item(140, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(140, 3) "Pieter De Lange"
Set-TextCell $ws.Cells.Item(140, 4) "Amaury Raynel"
Set-TextCell $ws.Cells.Item(140, 5) "Gana Pieter De Lange"
$ws.Cells.Item(140, 6).Value = 2.5

# row 141
$ws.Cells.Item(141, 1).Value = 14610356
Set-TextCell $ws.Cells.Item(141, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(141, 3) "Maik Steiner"
Set-TextCell $ws.Cells.Item(141, 4) "Stijn Paardekooper"
Set-TextCell $ws.Cells.Item(141, 5) "Gana Stijn Paardekooper"
$ws.Cells.Item(141, 6).Value = 2.5

# row 142
$ws.Cells.Item(142, 1).Value = 14610359
Set-TextCell $ws.Cells.Item(142, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(142, 3) "Mac Visser"
Set-TextCell $ws.Cells.Item(142, 4) "Nino Ehrenschneider"
Set-TextCell $ws.Cells.Item(142, 5) "Gana Mac Visser"
$ws.Cells.Item(142, 6).Value = 4.5

# row 143
$ws.Cells.Item(143, 1).Value = 14610360
Set-TextCell $ws.Cells.Item(143, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(143, 3) "Stijn Slump"
Set-TextCell $ws.Cells.Item(143, 4) "João Loureiro"
Set-TextCell $ws.Cells.Item(143, 5) "Gana João Loureiro"
$ws.Cells.Item(143, 6).Value = 4.5

# row 144
$ws.Cells.Item(144, 1).Value = 14609917
Set-TextCell $ws.Cells.Item(144, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(144, 3) "Amr Elsayed"
Set-TextCell $ws.Cells.Item(144, 4) "Edison Ambarzumjan"
Set-TextCell $ws.Cells.Item(144, 5) "Gana Edison Ambarzumjan"
$ws.Cells.Item(144, 6).Value = 3

# row 145
$ws.Cells.Item(145, 1).Value = 14609919
Set-TextCell $ws.Cells.Item(145, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(145, 3) "Kris van Wyk"
Set-TextCell $ws.Cells.Item(145, 4) "Karim Ibrahim"
Set-TextCell $ws.Cells.Item(145, 5) "Gana Karim Ibrahim"
$ws.Cells.Item(145, 6).Value = 4

# row 146
$ws.Cells.Item(146, 1).Value = 14609923
Set-TextCell $ws.Cells.Item(146, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(146, 3) "David Eichenseher"
Set-TextCell $ws.Cells.Item(146, 4) "Lorenzo Lorusso"
Set-TextCell $ws.Cells.Item(146, 5) "Gana David Eichenseher"
$ws.Cells.Item(146, 6).Value = 3.4

# row 147
$ws.Cells.Item(147, 1).Value = 14609920
Set-TextCell $ws.Cells.Item(147, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(147, 3) "Ivan Gretskiy"
Set-TextCell $ws.Cells.Item(147, 4) "Mohamed Safwat"
Set-TextCell $ws.Cells.Item(147, 5) "Gana Ivan Gretskiy"
$ws.Cells.Item(147, 6).Value = 5

# row 148
$ws.Cells.Item(148, 1).Value = 14609486
Set-TextCell $ws.Cells.Item(148, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(148, 3) "Giuseppe La Vela"
Set-TextCell $ws.Cells.Item(148, 4) "Marko Maksimovic"
Set-TextCell $ws.Cells.Item(148, 5) "Gana Marko Maksimovic"
$ws.Cells.Item(148, 6).Value = 2.25

# row 149
$ws.Cells.Item(149, 1).Value = 14609493
Set-TextCell $ws.Cells.Item(149, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(149, 3) "Andrey Chepelev"
Set-TextCell $ws.Cells.Item(149, 4) "Marko Miladinovic"
Set-TextCell $ws.Cells.Item(149, 5) "Gana Marko Miladinovic"
$ws.Cells.Item(149, 6).Value = 3.75

# row 150
$ws.Cells.Item(150, 1).Value = 14609494
Set-TextCell $ws.Cells.Item(150, 2) "2025-09-04"
Set-TextCell $ws.Cells.Item(150, 3) "Ognjen Milić"
Set-TextCell $ws.Cells.Item(150, 4) "Denys Klok"
Set-TextCell $ws.Cells.Item(150, 5) "Gana Denys Klok"
$ws.Cells.Item(150, 6).Value = 2.5

Write-Output "tracker actualizado: filas 120,121 con resultado; filas 129-150 anadidas"
